$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Quyết định bổ nhiệm LT 5.0"
$ws.Range("C4").Value = "<p>Quyết định bổ nhiệm LT 5.0</p>"
$ws.Range("E4").Value = "Ban Đào tạo"
$ws.Range("F4").Value = "12/09/2022 11:05"
$ws.Range("G4").Value = "2hfsdjfhjkadfhads"

# Row 5
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Quyết định bổ nhiệm gen 5.0"
$ws.Range("C5").Value = "<p>Quyết định bổ nhiệm gen 5.2</p>"
$ws.Range("E5").Value = "Ban Điều hành"
$ws.Range("F5").Value = "12/09/2022 11:13"
$ws.Range("G5").Value = "dfdaffdaff"
